$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.116.25"
$ws.Range("E2").Value = "  +2.91%  "

$ws.Range("D3").Value = "2.956.15"
$ws.Range("E3").Value = "  +1.11%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").Value = "'149.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "2.953.31"
$ws.Range("E8").Value = "  +1.12%  "

$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("D10").Value = "'7.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.98%  "

$ws.Range("E11").Value = "  +6.76%  "

$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("E13").Value = "  +5.30%  "

$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("D16").Value = "3.446.88"
$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("D17").Value = "63.070.02"
$ws.Range("E17").Value = "  +2.99%  "

$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D19").Value = "2.970.98"
$ws.Range("E19").Value = "  +1.54%  "

$ws.Range("D20").Value = "'442.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  -0.63%  "

$ws.Range("D23").Value = "'7.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "

$ws.Range("D24").Value = "'11.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.35%  "

$ws.Range("E25").Value = "  -0.95%  "

$ws.Range("D26").Value = "'2.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.79%  "

$ws.Range("E27").Value = "  +0.42%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "'7.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.56%  "

$ws.Range("D30").Value = "'2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("E32").Value = "  +16.30%  "

$ws.Range("E33").Value = "  -0.40%  "

$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "'0.992"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.56%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.94%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'5.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.23%  "

$ws.Range("E39").Value = "  +2.61%  "

$ws.Range("D40").Value = "'49.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("E41").Value = "  -0.39%  "

$ws.Range("E42").Value = "  -4.31%  "

$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").Value = "'39.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.67%  "

$ws.Range("D45").Value = "2.703.89"
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").Value = "'135.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.53%  "

$ws.Range("D47").Value = "'0.0339"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.70%  "

$ws.Range("D48").Value = "'360.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.39%  "

$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("D51").Value = "'22.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.74%  "

